$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from "Inventory" to "Sheet1"
$ws.Name = "Sheet1"

# Combine delivery date and invoice total rows into a single row 2
$ws.Range("B2").Value = "Delivery Date: 2024-10-12"
$ws.Range("E2").Value = "Invoice Total: 500"

# Clear old rows 3, 4, 6, 7 (data now consolidated into row 2)
$ws.Range("A3:G4").Clear()
$ws.Range("A6:G7").Clear()

# Apply bold, centered, thin-bordered style to header row 1 (A1:G1)
$headerRange = $ws.Range("A1:G1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
